$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '71.031.66'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -2.13%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.947.27'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -2.53%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.17%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '538.81'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +3.30%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '147.56'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '3.943.71'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -2.41%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.687'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -5.61%  '
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +0.01%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.740'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -5.65%  '
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -6.91%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '55.54'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +14.37%  '
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -4.29%  '
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -5.21%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '4.571.75'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -2.74%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.946.65'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -2.81%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '13.89'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -3.20%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '20.51'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -3.77%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '1.16'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -5.03%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '70.857.99'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -2.21%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '421.26'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -7.35%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '3.60'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -0.86%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '97.59'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -7.57%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '4.25'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +5.57%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '14.41'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -5.33%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '11.34'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -1.06%  '
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +14.95%  '
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -4.07%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '5.90'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +1.04%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '36.45'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -4.77%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '7.86'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +17.16%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '51.19'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +20.60%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.132'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +0.52%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '13.32'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -3.26%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '675.10'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -0.62%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '65.43'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -3.27%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.440'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +1.24%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0₃0814'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -6.47%  '
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -3.78%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.38'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -2.92%  '
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -0.07%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.998'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -0.11%  '
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -4.18%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.17'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -0.98%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '10.11'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +2.83%  '
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -5.68%  '
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -1.81%  '
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -2.92%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '3.01'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -2.89%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '144.03'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -0.69%  '
